# Update the "想去人数" (want-to-go count) figures in the F column on the
# "展览" and "全部类型" worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    2  = 694
    3  = 11
    4  = 530
    9  = 3547
    10 = 4296
    11 = 9
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
